# Insert a new weekly record at row 62, pushing the existing rows 62:95
# down to 63:96 (dimension grows from A1:R95 to A1:R96).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(62).Insert()

$ws.Range("A62").Value = 2
$ws.Range("B62").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C62").Value = 'Coquimbo'
$ws.Range("D62").Value = 44510
$ws.Range("E62").Value = 4
$ws.Range("F62").Value = 100112024
$ws.Range("G62").Value = 'Choclo'
$ws.Range("H62").Value = 'Dulce o Americano'
$ws.Range("I62").Value = 'Primera'
$ws.Range("J62").Value = 1100
$ws.Range("K62").Value = 28000
$ws.Range("L62").Value = 30000
$ws.Range("M62").Value = 29000
$ws.Range("N62").Value = '$/malla 70 unidades'
$ws.Range("O62").Value = 'Provincia de Limarí'
$ws.Range("P62").Value = 414
$ws.Range("Q62").Value = 70
$ws.Range("R62").Value = 'Hortaliza'
